$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that contain data in rows 2-25 (A = case index, G = 1 fixed, H = empty)
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# New vm_pu values per row (row number => array aligned with $cols)
$data = @{
    2 = @(1.02, 1.024852614231394, 1.030275598867084, 1.04973034196552, 1.054357701271882, 1.032927172094084, 1.030024908841102, 1.033087081237461, 1.052486587030788, 1.057101133296821, 1.031487662357165)
    3 = @(1.02, 1.025821940348247, 1.031003769423294, 1.05083822330453, 1.055499693330672, 1.033133484924805, 1.030633216276745, 1.033623788034601, 1.053406010297923, 1.058055512254551, 1.03209683365914)
    4 = @(1.02, 1.026449024731784, 1.031474435580884, 1.05155586022906, 1.056239099691689, 1.033265086062988, 1.031026115864915, 1.03396989909974, 1.054001111976669, 1.058672936690368, 1.032490291209791)
    5 = @(1.02, 1.026712618548812, 1.031672180864613, 1.051857737596656, 1.056550056504243, 1.033319956558579, 1.03119111895166, 1.034115122443009, 1.054251334189243, 1.058932472476968, 1.032655528619845)
    6 = @(1.02, 1.026756875199182, 1.031705375941115, 1.051908434854586, 1.05660227399346, 1.033329142883358, 1.031218813601369, 1.034139489540115, 1.054293350030483, 1.058976047947797, 1.032683262599134)
    7 = @(1.02, 1.026452547013156, 1.031477078346644, 1.051559893210931, 1.056243254278153, 1.033265821032186, 1.031028321320311, 1.033971840689088, 1.054004455295227, 1.058676404736092, 1.032492499797187)
    8 = @(1.02, 1.025180230202041, 1.030521792320598, 1.05010459753543, 1.054743547343114, 1.032997289003116, 1.030230637147537, 1.033268706526686, 1.052797274517352, 1.057423695450985, 1.031693682821392)
    9 = @(1.02, 1.022937227275699, 1.02883459972755, 1.047546028751686, 1.052104404015282, 1.032509592063097, 1.028819553024303, 1.032020730513451, 1.050671396273953, 1.055215326862569, 1.03028059479685)
    10 = @(1.02, 1.021441219748342, 1.027707264731466, 1.045844244656177, 1.050347346474936, 1.032174731268811, 1.027875183704065, 1.031182756298828, 1.04925503572529, 1.053742454865798, 1.029334884363856)
    11 = @(1.02, 1.020793274114707, 1.027218524302133, 1.045108283350173, 1.049587081778867, 1.032027430118315, 1.027465401306779, 1.03081849047747, 1.048641946187672, 1.053104535948163, 1.028924520028564)
    12 = @(1.02, 1.020552573551271, 1.027036895431589, 1.044835052878601, 1.04930476839519, 1.031972370063938, 1.027313060343817, 1.030682973462053, 1.048414247862361, 1.052867560913704, 1.028771962723961)
    13 = @(1.02, 1.020604205757063, 1.027075859452163, 1.044893655486873, 1.049365321833058, 1.031984196277569, 1.027345743854339, 1.030712051928884, 1.048463088552591, 1.052918393931147, 1.028804692648815)
    14 = @(1.02, 1.020773378234986, 1.027203512612436, 1.045085695213846, 1.04956374397953, 1.032022885890334, 1.027452811395967, 1.030807292930246, 1.048623123953593, 1.053084947999158, 1.028911912238634)
    15 = @(1.02, 1.020877607708559, 1.027282152178462, 1.045204035549547, 1.049686009372899, 1.032046678019304, 1.027518762086536, 1.030865945914569, 1.048721731071524, 1.053187564332545, 1.028977956586748)
    16 = @(1.02, 1.021484218301087, 1.027739688313213, 1.045893107396755, 1.050397814362318, 1.032184458642467, 1.027902361429907, 1.031206901599213, 1.049295728737238, 1.053784788150613, 1.029362100685186)
    17 = @(1.02, 1.021864685065331, 1.028026529371584, 1.046325590741067, 1.050844458634284, 1.03227026812044, 1.028142751998979, 1.031420394944681, 1.049655836756432, 1.054159369093032, 1.02960283263644)
    18 = @(1.02, 1.022086589019465, 1.028193781261168, 1.046577939978029, 1.051105032038785, 1.032320096991528, 1.028282884342871, 1.031544785204856, 1.049865901361599, 1.054377840815042, 1.029743163984333)
    19 = @(1.02, 1.022162249840725, 1.028250800031906, 1.04666399965365, 1.05119388989468, 1.032337049624962, 1.028330651693882, 1.031587175851124, 1.0499375313096, 1.054452331486421, 1.029790999170462)
    20 = @(1.02, 1.021823866223789, 1.027995759997683, 1.04627918016526, 1.050796532413764, 1.032261084571323, 1.028116968990904, 1.03139750327506, 1.049617198538757, 1.05411918165592, 1.029577013013535)
    21 = @(1.02, 1.020723561838681, 1.027165924386734, 1.045029140484322, 1.049505311320178, 1.032011502308141, 1.027421286246893, 1.030779252690517, 1.048575996678542, 1.053035902628667, 1.028880342320233)
    22 = @(1.02, 1.020031613865649, 1.026643659747861, 1.044243991107207, 1.048693948143566, 1.031852579222538, 1.026983133593486, 1.030389304957013, 1.047921527679429, 1.052354665801984, 1.028441567439801)
    23 = @(1.02, 1.020398442124697, 1.02692057058646, 1.044660137886009, 1.049124021826224, 1.031937016900646, 1.027215477528532, 1.030596139968093, 1.048268457472646, 1.052715815308898, 1.028674241329885)
    24 = @(1.02, 1.021842310550152, 1.028009663531415, 1.046300150831501, 1.050818188044779, 1.032265234908853, 1.02812861948009, 1.031407847457841, 1.04963465742515, 1.054137340675082, 1.029588680047753)
    25 = @(1.02, 1.023517216667257, 1.029271230979475, 1.048206787437954, 1.052786267609744, 1.032637390161139, 1.029184996725037, 1.032344420686054, 1.051220829482306, 1.055786353614676, 1.030646557469556)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value = $vals[$i]
    }
}

Write-Output "Updated $($data.Keys.Count) rows x $($cols.Length) columns of vm_pu results (380 kV case)"
